# DYNGROWTH_holdings.xlsx refresh: updates the "as of" disclaimer date and
# refreshes the Weight (col D) / Percent Change (col E) figures for the
# holdings table on Sheet1 (rows 2-77), matching a new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected (no password needed to unprotect here); drop
# protection so the locked cells can be written, then restore it afterwards.
$ws.Unprotect()

# --- Update the "as of" date in the disclaimer text (shared string) ---
$disclaimerCell = $ws.Range("A80")
$disclaimerCell.Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-20 for illustrative purposes only and are subject to change."

# --- Update Weight (col D) / Percent Change (col E) values, rows 2-77 ---
$ws.Range("D2").Value = 0.07486016806290953
$ws.Range("D3").Value = 0.04513581819486788
$ws.Range("D4").Value = 0.03553796580036234
$ws.Range("D5").Value = 0.03508481738404525
$ws.Range("D6").Value = 0.032899469396495
$ws.Range("D7").Value = 0.03013049407363111
$ws.Range("D8").Value = 0.03029863598600139
$ws.Range("D9").Value = 0.02885134354406234
$ws.Range("E9").Value = -0.005822704586865179
$ws.Range("D10").Value = 0.02634402168176048
$ws.Range("D11").Value = 0.02685798738553065
$ws.Range("D12").Value = 0.02484008568778179
$ws.Range("D13").Value = 0.0243787885499432
$ws.Range("D14").Value = 0.0198997742033985
$ws.Range("D15").Value = 0.0183684108052743
$ws.Range("D16").Value = 0.02012535466503003
$ws.Range("D17").Value = 0.01808757803674095
$ws.Range("D18").Value = 0.01712115953922259
$ws.Range("D19").Value = 0.01552505355486362
$ws.Range("D20").Value = 0.01392134540957301
$ws.Range("D21").Value = 0.01522827115457174
$ws.Range("D22").Value = 0.01458819901652385
$ws.Range("D23").Value = 0.01370893208942437
$ws.Range("D24").Value = 0.01520938997055852
$ws.Range("D25").Value = 0.01439143720417564
$ws.Range("D26").Value = 0.01249616382785995
$ws.Range("D27").Value = 0.01242779406680158
$ws.Range("D28").Value = 0.01234866699958832
$ws.Range("D29").Value = 0.01221453606210498
$ws.Range("D30").Value = 0.01071596629937215
$ws.Range("D31").Value = 0.01234233186547862
$ws.Range("D32").Value = 0.01266251699648162
$ws.Range("D33").Value = 0.01107406286038063
$ws.Range("D34").Value = 0.01147136272187969
$ws.Range("D35").Value = 0.00973613222373915
$ws.Range("D36").Value = 0.0108045588022552
$ws.Range("D37").Value = 0.01101816461823625
$ws.Range("D38").Value = 0.01020776432545865
$ws.Range("D39").Value = 0.009224949010244609
$ws.Range("D40").Value = 0.00938971218442306
$ws.Range("D41").Value = 0.009025702831578872
$ws.Range("D42").Value = 0.009959228319053154
$ws.Range("D43").Value = 0.009982283238479814
$ws.Range("D44").Value = 0.009285170049781483
$ws.Range("D45").Value = 0.009312199955316189
$ws.Range("D46").Value = 0.009597504583221094
$ws.Range("D47").Value = 0.009303852484489296
$ws.Range("D48").Value = 0.007311191948526515
$ws.Range("D49").Value = 0.008557350093398509
$ws.Range("D50").Value = 0.008021025092770587
$ws.Range("D51").Value = 0.008023310709782712
$ws.Range("D52").Value = 0.007771495339838085
$ws.Range("D53").Value = 0.007330868129761336
$ws.Range("D54").Value = 0.007192041739937876
$ws.Range("D55").Value = 0.006692137549524912
$ws.Range("D56").Value = 0.00636792774508752
$ws.Range("D57").Value = 0.007081040252870731
$ws.Range("D58").Value = 0.006604041919905373
$ws.Range("D59").Value = 0.005908319976344858
$ws.Range("D60").Value = 0.006162222213996211
$ws.Range("D61").Value = 0.005599761679707893
$ws.Range("D62").Value = 0.005741469934459684
$ws.Range("D63").Value = 0.00526288160705111
$ws.Range("D64").Value = 0.004955217682288454
$ws.Range("D65").Value = 0.004754083385221395
$ws.Range("D66").Value = 0.004380037192454393
$ws.Range("D67").Value = 0.004586935219378118
$ws.Range("D68").Value = 0.003553786642570964
$ws.Range("D69").Value = 0.003927932209990666
$ws.Range("D70").Value = 0.003933050004604774
$ws.Range("D71").Value = 0.00312255033717447
$ws.Range("D72").Value = 0.002528588037979899
$ws.Range("D73").Value = 0.00250816654684982
$ws.Range("D74").Value = 0.002182167998663809
$ws.Range("D75").Value = 0.001978748084584624
$ws.Range("D76").Value = 0.001966525002302387
$ws.Range("E77").Value = -0.002459686622529378

# Restore sheet protection to its original state
$ws.Protect("", $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true)
